$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new rows above the current row 50 (this shifts the existing
# rows 50-59 down to become rows 56-65, preserving all their data/styles).
$ws.Range("A50:A55").EntireRow.Insert()

# Data for the 6 newly-inserted rows (new week of Chirimoya prices for
# "Provincia de Limarí", dated 2021-10-05 / Excel serial 44474).
# Columns: A Mercado ID, B Mercado, C Region, D Fecha, E Codreg, F Tipo,
# G Producto ID, H Producto, I Categoria ID, J Categoria, K Variedad,
# L Calidad, M Volumen, N Precio Minimo, O Precio Maximo,
# P Precio Promedio, Q Unidad, R Comuna/Provincia, S Precio Kilo, T Factor

$newRows = @(
    @{ Row=50; L="Cuarta";                  M=200; N=1000;  O=1000;  P=1000;  Q="`$/kilo (en caja de 15 kilos)"; S=1000; T=1 },
    @{ Row=51; L="Especial";                M=180; N=20000; O=20000; P=20000; Q="`$/bandeja 8 kilos";            S=2500; T=8 },
    @{ Row=52; L="Extra (doble especial)";  M=250; N=24000; O=24000; P=24000; Q="`$/bandeja 8 kilos";            S=3000; T=8 },
    @{ Row=53; L="Primera";                 M=200; N=16000; O=16000; P=16000; Q="`$/bandeja 8 kilos";            S=2000; T=8 },
    @{ Row=54; L="Segunda";                 M=200; N=12800; O=12800; P=12800; Q="`$/bandeja 8 kilos";            S=1600; T=8 },
    @{ Row=55; L="Tercera";                 M=180; N=1400;  O=1400;  P=1400;  Q="`$/kilo (en caja de 15 kilos)"; S=1400; T=1 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 9
    $ws.Cells.Item($row, 2).Value = "Vega Central Mapocho de Santiago"
    $ws.Cells.Item($row, 3).Value = "Metropolitana"
    $ws.Cells.Item($row, 4).Value = 44474
    $ws.Cells.Item($row, 5).Value = 13
    $ws.Cells.Item($row, 6).Value = "Fruta"
    $ws.Cells.Item($row, 7).Value = 100107
    $ws.Cells.Item($row, 8).Value = "Otros"
    $ws.Cells.Item($row, 9).Value = 100107002
    $ws.Cells.Item($row, 10).Value = "Chirimoya"
    $ws.Cells.Item($row, 11).Value = "Cultivar IV Región"
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = "Provincia de Limarí"
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
}
